# Commit: Adjusted Mechanism coordinates based from URL (BCEF -> BCFG and
# CDGI -> CDEI). Also implemented positionSolver from PMKS onto MATLAB.
#
# The regenerated RMSE value for this run (Sheet1!A1) replaces the prior
# result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = 3.0869380577906647
